$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 63 (pushes the old rows 63..106 down to 65..108)
$ws.Rows.Item(63).Insert()
$ws.Rows.Item(63).Insert()

# New row 63: Santina / Especial
$ws.Cells.Item(63, 1).Value  = 7
$ws.Cells.Item(63, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(63, 3).Value  = "Ñuble"
$ws.Cells.Item(63, 4).Value  = 44596
$ws.Cells.Item(63, 5).Value  = 16
$ws.Cells.Item(63, 6).Value  = "Fruta"
$ws.Cells.Item(63, 7).Value  = 100103
$ws.Cells.Item(63, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(63, 9).Value  = 100103001
$ws.Cells.Item(63, 10).Value = "Cereza"
$ws.Cells.Item(63, 11).Value = "Santina"
$ws.Cells.Item(63, 12).Value = "Especial"
$ws.Cells.Item(63, 13).Value = 100
$ws.Cells.Item(63, 14).Value = 7000
$ws.Cells.Item(63, 15).Value = 7000
$ws.Cells.Item(63, 16).Value = 7000
$ws.Cells.Item(63, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(63, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(63, 19).Value = 700
$ws.Cells.Item(63, 20).Value = 10

# New row 64: Santina / Primera
$ws.Cells.Item(64, 1).Value  = 7
$ws.Cells.Item(64, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(64, 3).Value  = "Ñuble"
$ws.Cells.Item(64, 4).Value  = 44596
$ws.Cells.Item(64, 5).Value  = 16
$ws.Cells.Item(64, 6).Value  = "Fruta"
$ws.Cells.Item(64, 7).Value  = 100103
$ws.Cells.Item(64, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(64, 9).Value  = 100103001
$ws.Cells.Item(64, 10).Value = "Cereza"
$ws.Cells.Item(64, 11).Value = "Santina"
$ws.Cells.Item(64, 12).Value = "Primera"
$ws.Cells.Item(64, 13).Value = 200
$ws.Cells.Item(64, 14).Value = 5500
$ws.Cells.Item(64, 15).Value = 6000
$ws.Cells.Item(64, 16).Value = 5750
$ws.Cells.Item(64, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(64, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(64, 19).Value = 575
$ws.Cells.Item(64, 20).Value = 10
